$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns E/F/G ---
$ws.Range("E1").Value = "end"
$ws.Range("F1").Value = "days"
$ws.Range("G1").Value = "duration"

# --- Row 2 ---
# Move "days" value from G2 into F2
$ws.Range("F2").Value = "MWF"
# New "end" time for row 2 (2:00 PM), formatted like the start-time column
$ws.Range("E2").Value = 0.58333333333333337
$ws.Range("E2").NumberFormat = $ws.Range("D2").NumberFormat
# New "duration" formula column
$ws.Range("G2").Formula = "=E2-D2"
$ws.Range("G2").NumberFormat = $ws.Range("D2").NumberFormat

# --- Row 3 ---
$ws.Range("F3").Value = "TR"
$ws.Range("E3").Value = 0.70833333333333337
$ws.Range("E3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("G3").Formula = "=E3-D3"
$ws.Range("G3").NumberFormat = $ws.Range("D2").NumberFormat

# --- Extend the "duration" formatting down through row 7 (blank, formatted cells) ---
$ws.Range("G4:G7").NumberFormat = $ws.Range("D2").NumberFormat

# --- Selection cosmetics ---
$ws.Range("I6").Select()
